# Update TPM-derived numeric values for rows 2-9, then remove the now
# obsolete "Resolving-Mac"-as-sender rows (originally rows 10-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> Il13/Il4ra -> ECs)
$ws.Range("I2").Value = 0.7099439172299504
$ws.Range("J2").Value = 0.7099439172299504
$ws.Range("M2").Value = 26.93692533333333
$ws.Range("N2").Value = 80.810776
$ws.Range("O2").Value = 0.319184224619106
$ws.Range("P2").Value = 0.319184224619106
$ws.Range("Q2").Value = 11.46297265969956
$ws.Range("R2").Value = 103.166753937296
$ws.Range("S2").Value = 0.2266028987440925
$ws.Range("T2").Value = 0.2266028987440925

# Row 3 (FAPs -> Il13/Il4ra -> FAPs)
$ws.Range("I3").Value = 0.7099439172299504
$ws.Range("J3").Value = 0.7099439172299504
$ws.Range("O3").Value = 0.2530310992806105
$ws.Range("P3").Value = 0.2530310992806105
$ws.Range("S3").Value = 0.1796378898042771
$ws.Range("T3").Value = 0.1796378898042771

# Row 4 (FAPs -> Il13/Il4ra -> MuSCs)
$ws.Range("I4").Value = 0.7099439172299504
$ws.Range("J4").Value = 0.7099439172299504
$ws.Range("M4").Value = 9.099134999999999
$ws.Range("N4").Value = 27.297405
$ws.Range("O4").Value = 0.1078185544096088
$ws.Range("P4").Value = 0.1078185544096088
$ws.Range("Q4").Value = 3.872124767069999
$ws.Range("R4").Value = 34.84912290363
$ws.Range("S4").Value = 0.07654512686762818
$ws.Range("T4").Value = 0.0765451268676282

# Row 5 (FAPs -> Il13/Il4ra -> Resolving-Mac)
$ws.Range("I5").Value = 0.7099439172299504
$ws.Range("J5").Value = 0.7099439172299504
$ws.Range("M5").Value = 27.002912
$ws.Range("N5").Value = 81.008736
$ws.Range("O5").Value = 0.3199661216906747
$ws.Range("P5").Value = 0.3199661216906747
$ws.Range("Q5").Value = 11.49105319771733
$ws.Range("R5").Value = 103.419478779456
$ws.Range("S5").Value = 0.2271580018139526
$ws.Range("T5").Value = 0.2271580018139526

# Row 6 (MuSCs -> Il13/Il4ra -> ECs)
$ws.Range("I6").Value = 0.2900560827700495
$ws.Range("J6").Value = 0.2900560827700495
$ws.Range("M6").Value = 26.93692533333333
$ws.Range("N6").Value = 80.810776
$ws.Range("O6").Value = 0.319184224619106
$ws.Range("P6").Value = 0.319184224619106
$ws.Range("Q6").Value = 4.683334649229334
$ws.Range("R6").Value = 42.150011843064
$ws.Range("S6").Value = 0.09258132587501348
$ws.Range("T6").Value = 0.09258132587501348

# Row 7 (MuSCs -> Il13/Il4ra -> FAPs)
$ws.Range("I7").Value = 0.2900560827700495
$ws.Range("J7").Value = 0.2900560827700495
$ws.Range("O7").Value = 0.2530310992806105
$ws.Range("P7").Value = 0.2530310992806105
$ws.Range("S7").Value = 0.07339320947633336
$ws.Range("T7").Value = 0.07339320947633336

# Row 8 (MuSCs -> Il13/Il4ra -> MuSCs)
$ws.Range("I8").Value = 0.2900560827700495
$ws.Range("J8").Value = 0.2900560827700495
$ws.Range("M8").Value = 9.099134999999999
$ws.Range("N8").Value = 27.297405
$ws.Range("O8").Value = 0.1078185544096088
$ws.Range("P8").Value = 0.1078185544096088
$ws.Range("Q8").Value = 1.582002908505
$ws.Range("R8").Value = 14.238026176545
$ws.Range("S8").Value = 0.03127342754198056
$ws.Range("T8").Value = 0.03127342754198057

# Row 9 (MuSCs -> Il13/Il4ra -> Resolving-Mac)
$ws.Range("I9").Value = 0.2900560827700495
$ws.Range("J9").Value = 0.2900560827700495
$ws.Range("M9").Value = 27.002912
$ws.Range("N9").Value = 81.008736
$ws.Range("O9").Value = 0.3199661216906747
$ws.Range("P9").Value = 0.3199661216906747
$ws.Range("Q9").Value = 4.694807289056
$ws.Range("R9").Value = 42.253265601504
$ws.Range("S9").Value = 0.09280811987672209
$ws.Range("T9").Value = 0.09280811987672209

# The Resolving-Mac-as-sender block (former rows 10-13) is no longer part
# of this LR-pair's data; remove those trailing rows entirely.
$ws.Rows("10:13").Delete()
